$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Problem Statement Approved" column (J) as "Yes" for the
# rows that have since had their problem statement approved:
#   row 5  -> Guttapati,Devi Prasad Reddy
#   row 7  -> Ma,Xiaoye (Marshall)
#   row 10 -> Schaap,Alexander
$ws.Range("J5").Value = "Yes"
$ws.Range("J7").Value = "Yes"
$ws.Range("J10").Value = "Yes"

# Update the active selection on the sheet to J11
$ws.Range("J11").Select() | Out-Null
